# Update the "想去人数" (attendance interest count) figures on both the
# "展览" and "全部类型" worksheets, which hold duplicate data tables.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1386
    3 = 2210
    4 = 345
    6 = 6420
    7 = 289
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
